$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells retain their original text formatting (avoid Excel auto-converting
# numeric-looking strings like "0.7126" or "1.001" into actual numbers).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.343.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.53%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.58%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7126"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.62%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.75"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.28%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07796"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.76%  "

# Row 9
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3113"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.87%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.10"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.53%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08420"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.29%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.881.93"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.92%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.91%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7122"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.08"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.351.91"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.44%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.70%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008255"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.57%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.02"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.08%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.20"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.82%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.121.99"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.25%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.749"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.69%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1586"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.37%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.09"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.044"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.33%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.30%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.510"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.76%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.419"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.29%  "

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.289"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.60%  "

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.320"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.37%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05287"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.11%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.936"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.19%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.178"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.22%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7401"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -8.49%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.699"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.69%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.85%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.214.47"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.89%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.729"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.543"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.86%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.89"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +8.56%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8873"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.29%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.89"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.04%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.017.79"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.19%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.803"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.90%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5212"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.55%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000122"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.90%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.403"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.41%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.04%  "
